$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Cell="D2"; Text="303.97"},
    @{Cell="E2"; Text="-1.62%"},
    @{Cell="G2"; Text="6"},
    @{Cell="D3"; Text="35.70"},
    @{Cell="E3"; Text="-1.34%"},
    @{Cell="G3"; Text="6"},
    @{Cell="D4"; Text="5.041"},
    @{Cell="E4"; Text="-1.47%"},
    @{Cell="G4"; Text="6"},
    @{Cell="D5"; Text="0.07951"},
    @{Cell="E5"; Text="-2.25%"},
    @{Cell="G5"; Text="6"},
    @{Cell="D6"; Text="1.873"},
    @{Cell="E6"; Text="-4.87%"},
    @{Cell="G6"; Text="6"},
    @{Cell="B7"; Text="GateToken"},
    @{Cell="C7"; Text="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"},
    @{Cell="D7"; Text="4.104"},
    @{Cell="E7"; Text="-1.92%"},
    @{Cell="G7"; Text="6"},
    @{Cell="B8"; Text="KuCoinToken"},
    @{Cell="C8"; Text="https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"},
    @{Cell="D8"; Text="7.747"},
    @{Cell="E8"; Text="-0.53%"},
    @{Cell="G8"; Text="6"},
    @{Cell="B9"; Text="MXToken"},
    @{Cell="C9"; Text="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"},
    @{Cell="D9"; Text="0.9202"},
    @{Cell="E9"; Text="-1.02%"},
    @{Cell="G9"; Text="6"},
    @{Cell="B10"; Text="LiechtensteinCryptoassetsExchange"},
    @{Cell="C10"; Text="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"},
    @{Cell="D10"; Text="0.1347"},
    @{Cell="E10"; Text="-1.99%"},
    @{Cell="G10"; Text="6"},
    @{Cell="B11"; Text="WazirX"},
    @{Cell="C11"; Text="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"},
    @{Cell="D11"; Text="0.1883"},
    @{Cell="E11"; Text="-2.34%"},
    @{Cell="G11"; Text="6"},
    @{Cell="B12"; Text="MandalaExchangeToken"},
    @{Cell="C12"; Text="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"},
    @{Cell="D12"; Text="0.09034"},
    @{Cell="E12"; Text="-2.35%"},
    @{Cell="G12"; Text="6"},
    @{Cell="B13"; Text="BitrueCoin"},
    @{Cell="C13"; Text="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"},
    @{Cell="D13"; Text="0.03440"},
    @{Cell="E13"; Text="0.96%"},
    @{Cell="G13"; Text="6"},
    @{Cell="B14"; Text="BitMartToken"},
    @{Cell="C14"; Text="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"},
    @{Cell="D14"; Text="0.09819"},
    @{Cell="E14"; Text="-0.23%"},
    @{Cell="G14"; Text="6"},
    @{Cell="B15"; Text="BitForexToken"},
    @{Cell="C15"; Text="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"},
    @{Cell="D15"; Text="0.001411"},
    @{Cell="E15"; Text="-0.29%"},
    @{Cell="G15"; Text="6"},
    @{Cell="B16"; Text="TigerCash"},
    @{Cell="C16"; Text="https://coinranking.com/coin/6hIn06L2+tigercash-tch"},
    @{Cell="D16"; Text="0.006012"},
    @{Cell="E16"; Text="4.53%"},
    @{Cell="G16"; Text="6"},
    @{Cell="B17"; Text="LEO"},
    @{Cell="C17"; Text="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"},
    @{Cell="D17"; Text="3.736"},
    @{Cell="E17"; Text="3.14%"},
    @{Cell="G17"; Text="6"},
    @{Cell="D18"; Text="3.393"},
    @{Cell="E18"; Text="14.24%"},
    @{Cell="G18"; Text="6"},
    @{Cell="E19"; Text="0.15%"},
    @{Cell="G19"; Text="6"},
    @{Cell="D20"; Text="0.1333"},
    @{Cell="E20"; Text="2.21%"},
    @{Cell="G20"; Text="6"},
    @{Cell="D21"; Text="5.176"},
    @{Cell="E21"; Text="5.91%"},
    @{Cell="G21"; Text="6"},
    @{Cell="D22"; Text="0.2348"},
    @{Cell="E22"; Text="-6.02%"},
    @{Cell="G22"; Text="6"},
    @{Cell="D23"; Text="0.04403"},
    @{Cell="E23"; Text="-1.17%"},
    @{Cell="G23"; Text="6"},
    @{Cell="D24"; Text="0.001196"},
    @{Cell="E24"; Text="-1.80%"},
    @{Cell="G24"; Text="6"},
    @{Cell="D25"; Text="0.004618"},
    @{Cell="E25"; Text="-5.28%"},
    @{Cell="G25"; Text="6"},
    @{Cell="D26"; Text="0.0001297"},
    @{Cell="E26"; Text="4.46%"},
    @{Cell="G26"; Text="6"},
    @{Cell="D27"; Text="0.0004437"},
    @{Cell="E27"; Text="-0.09%"},
    @{Cell="G27"; Text="6"},
    @{Cell="G28"; Text="6"},
    @{Cell="G29"; Text="6"},
    @{Cell="G30"; Text="6"},
    @{Cell="G31"; Text="6"},
    @{Cell="G32"; Text="6"},
    @{Cell="G33"; Text="6"},
    @{Cell="G34"; Text="6"},
    @{Cell="G35"; Text="6"},
    @{Cell="G36"; Text="6"},
    @{Cell="G37"; Text="6"},
    @{Cell="G38"; Text="6"},
    @{Cell="E39"; Text="-4.77%"},
    @{Cell="G39"; Text="6"},
    @{Cell="D40"; Text="0.05265"},
    @{Cell="E40"; Text="5.85%"},
    @{Cell="G40"; Text="6"},
    @{Cell="E41"; Text="-0.66%"},
    @{Cell="G41"; Text="6"},
    @{Cell="E42"; Text="-1.42%"},
    @{Cell="G42"; Text="6"},
    @{Cell="E43"; Text="-2.85%"},
    @{Cell="G43"; Text="6"},
    @{Cell="D44"; Text="0.002145"},
    @{Cell="E44"; Text="2.00%"},
    @{Cell="G44"; Text="6"},
    @{Cell="D45"; Text="0.01016"},
    @{Cell="E45"; Text="-15.04%"},
    @{Cell="G45"; Text="6"},
    @{Cell="D46"; Text="0.00006156"},
    @{Cell="E46"; Text="-4.44%"},
    @{Cell="G46"; Text="6"},
    @{Cell="E47"; Text="-0.10%"},
    @{Cell="G47"; Text="6"},
    @{Cell="G48"; Text="6"},
    @{Cell="E49"; Text="39.21%"},
    @{Cell="G49"; Text="6"},
    @{Cell="E50"; Text="-0.10%"},
    @{Cell="G50"; Text="6"},
    @{Cell="E51"; Text="-0.10%"},
    @{Cell="G51"; Text="6"}
)

foreach ($e in $edits) {
    $r = $ws.Range($e.Cell)
    $r.NumberFormat = "@"
    $r.Value = $e.Text
    $r.Style = "Normal"
}
